$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 339
$ws1.Range("F3").Value = 274
$ws1.Range("F4").Value = 32
$ws1.Range("F5").Value = 3349
$ws1.Range("F6").Value = 2134
$ws1.Range("F11").Value = 1241
$ws1.Range("F12").Value = 224
$ws1.Range("F13").Value = 1392

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 339
$ws4.Range("F3").Value = 274
$ws4.Range("F4").Value = 32
$ws4.Range("F5").Value = 3349
$ws4.Range("F6").Value = 2134
$ws4.Range("F14").Value = 1241
$ws4.Range("F15").Value = 224
$ws4.Range("F16").Value = 1392
